# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with per-fund holding data (same shape as the other
#    quarterly sheets).
# 2. Prepend a new row to the "总计" sheet summarising the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Locate the existing "2021-Q4" sheet.
# ---------------------------------------------------------------------
$q4sheet  = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned between "2021-Q4" and "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $q4sheet)
$newSheet.Name = "2022-Q1"

# Re-fetch "总计" by name now that the sheet collection has shifted — a
# reference obtained before the insert would otherwise track the (now
# stale) positional slot rather than the worksheet itself.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (row 1), columns B..H — bold, thin border, centered/top.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows (row 2..6): A=index(number), B..G=text, H=rank(number).
$rows = @(
    @("004932", "招商丰拓灵活配置混合A", "24.42", "46.69", "1.84", "0.4493", 9),
    @("580008", "东吴新产业精选股票A", "2.32", "89.77", "9.73", "0.2257", 1),
    @("004933", "招商丰拓灵活配置混合C", "10.85", "46.69", "1.84", "0.1996", 9),
    @("010740", "汇安核心价值混合A", "0.95", "93.68", "3.39", "0.0322", 5),
    @("010741", "汇安核心价值混合C", "0.36", "93.68", "3.39", "0.0122", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $idxCell = $newSheet.Cells.Item($excelRow, 1)
    $idxCell.Value = $r
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    for ($c = 0; $c -lt 6; $c++) {
        $cell = $newSheet.Cells.Item($excelRow, 2 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
    }

    $rankCell = $newSheet.Cells.Item($excelRow, 8)
    $rankCell.Value = $row[6]
}

# ---------------------------------------------------------------------
# 2. Insert a new row into the "总计" sheet for 2022-Q1 and renumber the
#    index column. Existing rows 2..6 are pushed down to 3..7 by copying
#    values (bottom-up, so nothing is clobbered before it is read) rather
#    than using Rows.Insert(), which would drag the header row's bold /
#    bordered formatting onto the shifted data cells.
# ---------------------------------------------------------------------
for ($r = 6; $r -ge 2; $r--) {
    $dest = $r + 1
    $totalSheet.Cells.Item($dest, 2).Value = $totalSheet.Cells.Item($r, 2).Value2
    $totalSheet.Cells.Item($dest, 3).Value = $totalSheet.Cells.Item($r, 3).Value2
    $totalSheet.Cells.Item($dest, 4).Value = $totalSheet.Cells.Item($r, 4).Value2
}

# Row 7's index cell (A7) is brand new — give it the same bold/border/
# center-top styling as the rest of the index column (A2:A6).
$a7 = $totalSheet.Cells.Item(7, 1)
$a7.Font.Bold = $true
$a7.HorizontalAlignment = -4108
$a7.VerticalAlignment = -4160
$a7.Borders.LineStyle = 1

# Renumber the index column (A2:A7) to 0..5.
for ($r = 2; $r -le 7; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# Write the new 2022-Q1 summary row.
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 5
$totalSheet.Cells.Item(2, 4).Value = 0.92

Write-Host "2022-Q1 sheet added and 总计 updated"
